$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price (column D) values must stay text, matching the
# original inline-string cell type -- force Text format before assigning,
# then drop back to the workbook default style so no stray formatting is
# left behind on the cell.
$textCells = @("D5", "D6", "D10", "D13", "D16", "D18", "D19", "D21", "D22", "D23", "D24", "D26", "D28", "D29", "D32", "D34", "D35", "D37", "D38", "D39", "D40", "D43", "D45", "D46", "D48", "D50")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.081.00"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.599.35"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "583.49"
$ws.Range("E5").Value = "  +2.53%  "
$ws.Range("D6").Value = "147.87"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  +2.81%  "
$ws.Range("E9").Value = "  +3.01%  "
$ws.Range("D10").Value = "5.68"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("E12").Value = "  +0.10%  "
$ws.Range("D13").Value = "27.25"
$ws.Range("E13").Value = "  +0.23%  "
$ws.Range("D14").Value = "3.066.03"
$ws.Range("E14").Value = "  +2.10%  "
$ws.Range("D15").Value = "62.998.10"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "0.0000147"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").Value = "2.601.33"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("D18").Value = "11.36"
$ws.Range("E18").Value = "  +0.53%  "
$ws.Range("D19").Value = "343.28"
$ws.Range("E19").Value = "  +2.58%  "
$ws.Range("E20").Value = "  +1.41%  "
$ws.Range("D21").Value = "6.78"
$ws.Range("E21").Value = "  +0.28%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "5.72"
$ws.Range("E23").Value = "  -0.56%  "
$ws.Range("D24").Value = "67.15"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "2.725.90"
$ws.Range("E25").Value = "  +1.97%  "
$ws.Range("D26").Value = "0.169"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("E27").Value = "  +0.68%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "8.40"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  +7.69%  "
$ws.Range("E31").Value = "  -0.58%  "
$ws.Range("D32").Value = "1.94"
$ws.Range("E32").Value = "  +5.30%  "
$ws.Range("D33").Value = "0.0₃0822"
$ws.Range("E33").Value = "  +1.32%  "
$ws.Range("D34").Value = "465.87"
$ws.Range("E34").Value = "  +15.08%  "
$ws.Range("D35").Value = "176.83"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("E36").Value = "  +4.28%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("B38").Value = "PolygonEcosystemToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D38").Value = "0.402"
$ws.Range("E38").Value = "  +0.51%  "
$ws.Range("D39").Value = "19.23"
$ws.Range("E39").Value = "  +0.71%  "
$ws.Range("D40").Value = "4.59"
$ws.Range("E40").Value = "  +5.84%  "
$ws.Range("D43").Value = "159.24"
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("E44").Value = "  +1.30%  "
$ws.Range("D45").Value = "0.639"
$ws.Range("E45").Value = "  +6.45%  "
$ws.Range("D46").Value = "21.05"
$ws.Range("E46").Value = "  +1.62%  "
$ws.Range("E47").Value = "  +3.13%  "
$ws.Range("D48").Value = "0.0974"
$ws.Range("E48").Value = "  +1.00%  "
$ws.Range("E49").Value = "  -0.69%  "
$ws.Range("D50").Value = "18.59"
$ws.Range("E50").Value = "  +1.71%  "
$ws.Range("E51").Value = "  -0.18%  "

foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
